# Update the Des Moines Antibiogram worksheet with the refreshed
# Jan-Dec 2024 susceptibility data and refresh the reporting-period /
# revision-date labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Reporting period label (row 3) ---------------------------------
$ws.Range("F3").Value = "Jan. 1st, 2024 to Dec. 31st, 2024"

# --- Gram-negative panel (rows 11-15) --------------------------------
# E.coli
$ws.Range("C11").Value = 167
$ws.Range("D11").Value = 98
$ws.Range("E11").Value = 57
$ws.Range("F11").Value = 96
$ws.Range("G11").Value = 99
$ws.Range("H11").Value = 92
$ws.Range("I11").Value = 75
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 94
$ws.Range("L11").Value = 100
$ws.Range("M11").Value = 75
$ws.Range("N11").Value = 95
$ws.Range("O11").Value = 93
$ws.Range("P11").Value = 93
$ws.Range("Q11").Value = 74

# Pseud. aeruginosa
$ws.Range("C12").Value = 63
$ws.Range("D12").Value = 100
$ws.Range("F12").Value = 92
$ws.Range("G12").Value = 97
$ws.Range("J12").Value = 97
$ws.Range("K12").Value = 90
$ws.Range("L12").Value = 93
$ws.Range("M12").Value = 92
$ws.Range("N12").Value = 90

# Klebsiella pneumoniae
$ws.Range("C13").Value = 63
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 98
$ws.Range("G13").Value = 100
$ws.Range("H13").Value = 95
$ws.Range("I13").Value = 79
$ws.Range("J13").Value = 90
$ws.Range("K13").Value = 92
$ws.Range("N13").Value = 92
$ws.Range("O13").Value = 95
$ws.Range("P13").Value = 86
$ws.Range("Q13").Value = 83

# Proteus mirabilis
$ws.Range("C14").Value = 36
$ws.Range("E14").Value = 83
$ws.Range("F14").Value = 100
$ws.Range("G14").Value = 100
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 67
$ws.Range("J14").Value = 83
$ws.Range("K14").Value = 89
$ws.Range("L14").Value = 8
$ws.Range("M14").Value = 86
$ws.Range("O14").Value = 89
$ws.Range("P14").Value = 83
$ws.Range("Q14").Value = 89

# Enterobacter cloacae cmp.
$ws.Range("C15").Value = 34
$ws.Range("F15").Value = 94
$ws.Range("G15").Value = 99
$ws.Range("H15").Value = 88
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 97
$ws.Range("K15").Value = 95
$ws.Range("L15").Value = 94
$ws.Range("M15").Value = 94
$ws.Range("N15").Value = 95

# --- Gram-positive panel (rows 22-25) --------------------------------
# Staph aureus
$ws.Range("C22").Value = 74
$ws.Range("E22").Value = 81
$ws.Range("G22").Value = 59
$ws.Range("I22").Value = 92
$ws.Range("N22").Value = 99
$ws.Range("O22").Value = 86
$ws.Range("P22").Value = 99

# Staph epidermidis
$ws.Range("C23").Value = 53
$ws.Range("F23").Value = 71
$ws.Range("G23").Value = 49
$ws.Range("H23").Value = 97
$ws.Range("I23").Value = 86
$ws.Range("L23").Value = 32
$ws.Range("N23").Value = 94
$ws.Range("O23").Value = 81
$ws.Range("P23").Value = 72

# Entero faecalis
$ws.Range("C24").Value = 63
$ws.Range("E24").Value = 100
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 100
$ws.Range("O24").Value = 37

# Strep pneum
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 92
$ws.Range("K25").Value = 90
$ws.Range("M25").Value = 75

# --- Combined total patient count (row 28) ---------------------------
$ws.Range("C28").Value = 488

# --- Footer revision date (row 47) ------------------------------------
$ws.Range("A47").Value = "February 2025"

# --- Leave the selection where the author last left it -----------------
$ws.Range("A47").Select()
